$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06541609988146302
$ws.Range("H2").Value = 1.732288475226841
$ws.Range("I2").Value = -12.84803484966236
$ws.Range("G3").Value = 0.06595803288105929
$ws.Range("H3").Value = 17.34834830656951
$ws.Range("G4").Value = -0.01630143893053489
$ws.Range("H4").Value = 19.7432379859169
$ws.Range("G5").Value = -0.01737438082549005
$ws.Range("H5").Value = -51.98913302203832
$ws.Range("G6").Value = -0.006857526551817841
$ws.Range("H6").Value = 38.73684529492827
$ws.Range("G7").Value = -0.006564214915726695
$ws.Range("H7").Value = -15.899026768636
$ws.Range("G8").Value = 0.003049534851035189
$ws.Range("H8").Value = 153.6025220676941
$ws.Range("G9").Value = 0.006552494483497295
$ws.Range("H9").Value = 219.284908870346
$ws.Range("G10").Value = -0.05976204688561625
$ws.Range("H10").Value = 5.172353237690848
$ws.Range("G11").Value = -0.06801204634979652
$ws.Range("H11").Value = -6.139385461277814
$ws.Range("G12").Value = -0.383875069823512
$ws.Range("H12").Value = 2.718762940722252
$ws.Range("G13").Value = -0.4163681986694099
$ws.Range("H13").Value = -6.184088542560611
$ws.Range("G14").Value = -0.02272436342732542
$ws.Range("H14").Value = -177.3203886351492
$ws.Range("G15").Value = -0.01173063463570371
$ws.Range("H15").Value = 74.11523060403842
$ws.Range("G16").Value = 0.1443419171300175
$ws.Range("H16").Value = 5.551946809595822
$ws.Range("G17").Value = 0.1476464535448672
$ws.Range("H17").Value = 5.862033124567299
$ws.Range("G18").Value = 0.1265153754858295
$ws.Range("H18").Value = 7.51874089585674
$ws.Range("G19").Value = 0.1171332704777653
$ws.Range("H19").Value = -8.9653250552344
$ws.Range("G20").Value = 0.08761312955697763
$ws.Range("H20").Value = -1.26388435267985
$ws.Range("G21").Value = 0.0927195166229722
$ws.Range("H21").Value = 6.462768509479335
$ws.Range("G22").Value = -0.09904199692640561
$ws.Range("H22").Value = -5.946340453514313
$ws.Range("G23").Value = -0.1030400515177863
$ws.Range("H23").Value = -1.57231318300617
$ws.Range("G24").Value = 0.1627560847786358
$ws.Range("H24").Value = 1.033591358955108
$ws.Range("G25").Value = 0.1653490158232516
$ws.Range("H25").Value = -3.075526028323837
$ws.Range("G26").Value = 0.0831333107993223
$ws.Range("H26").Value = -8.295830122068709
$ws.Range("G27").Value = 0.08844278774595928
$ws.Range("H27").Value = 2.886215757295472
$ws.Range("G28").Value = -0.1377105337885386
$ws.Range("H28").Value = -0.05092793169942775
$ws.Range("G29").Value = -0.1405649722066565
$ws.Range("H29").Value = -0.5473590998903427
$ws.Range("G30").Value = 0.0504747709955134
$ws.Range("H30").Value = -2.961636954169125
$ws.Range("G31").Value = 0.03747939438275949
$ws.Range("H31").Value = -14.46849613243487
$ws.Range("G32").Value = 0.1140968152882486
$ws.Range("H32").Value = 4.961075041485894
$ws.Range("G33").Value = 0.1208720615263008
$ws.Range("H33").Value = -2.597305220494612
$ws.Range("G34").Value = -0.01532417659944549
$ws.Range("H34").Value = 1.87122166951309
$ws.Range("G35").Value = -0.0144018616024453
$ws.Range("H35").Value = 13.95343587258256
$ws.Range("G36").Value = 0.03123225943698414
$ws.Range("H36").Value = -15.05334594844025
$ws.Range("G37").Value = 0.03259832273360197
$ws.Range("H37").Value = -8.651449147130332
$ws.Range("G38").Value = 0.1042096727759903
$ws.Range("H38").Value = 3.895669543362004
$ws.Range("G39").Value = 0.1040239924792556
$ws.Range("H39").Value = 6.800797220294746
$ws.Range("G40").Value = 0.03419821709323458
$ws.Range("H40").Value = 1.515682467302804
$ws.Range("G41").Value = 0.03223725225220001
$ws.Range("H41").Value = 0.05750450344260649
$ws.Range("G42").Value = 0.124383150708606
$ws.Range("H42").Value = 2.876397039592414
$ws.Range("G43").Value = 0.117890355878734
$ws.Range("H43").Value = -7.743519741487746
$ws.Range("G44").Value = 0.04050944153364083
$ws.Range("H44").Value = 2.124713194916691
$ws.Range("G45").Value = 0.04050944153364083
$ws.Range("H45").Value = 29.99481015783694
$ws.Range("G46").Value = 0.06588728875228977
$ws.Range("H46").Value = 16.37231726129659
$ws.Range("G47").Value = 0.06632527107476163
$ws.Range("H47").Value = 13.04765648035194
$ws.Range("G48").Value = 0.0483458036386855
$ws.Range("H48").Value = -1.850525025195072
$ws.Range("G49").Value = 0.04402552233612568
$ws.Range("H49").Value = -3.384466067962066
$ws.Range("G50").Value = 0.0276239437082732
$ws.Range("H50").Value = 4.29449210892596
$ws.Range("G51").Value = 0.03442717896560511
$ws.Range("H51").Value = 22.88456813238026
$ws.Range("G52").Value = -0.08735847597933077
$ws.Range("H52").Value = -0.5002492562916754
$ws.Range("G53").Value = -0.08233357014995947
$ws.Range("H53").Value = -2.639576756276863
$ws.Range("G54").Value = 0.04594042701214705
$ws.Range("H54").Value = -8.168187703260019
$ws.Range("G55").Value = 0.04692816982865802
$ws.Range("H55").Value = -16.6241050192608
$ws.Range("G56").Value = 0.04674062402866477
$ws.Range("H56").Value = -5.444401756004499
$ws.Range("G57").Value = 0.04196317690557763
$ws.Range("H57").Value = 10.47599348937258
$ws.Range("G58").Value = 0.0522427165982714
$ws.Range("H58").Value = -9.322961618156288
$ws.Range("G59").Value = 0.06170559162566217
$ws.Range("H59").Value = 8.230152577872827
$ws.Range("G60").Value = 0.0266054934979764
$ws.Range("H60").Value = -3.111820959424944
$ws.Range("G61").Value = 0.03442279123370657
$ws.Range("H61").Value = 28.93481332359972
$ws.Range("G62").Value = 0.06220795609703562
$ws.Range("H62").Value = -0.3831637848147654
$ws.Range("G63").Value = 0.06630011994870408
$ws.Range("H63").Value = 3.784204102542447
$ws.Range("G64").Value = 0.03131818373618619
$ws.Range("H64").Value = 12.90039836423029
$ws.Range("G65").Value = 0.03149116796827313
$ws.Range("H65").Value = -11.11032815647788
$ws.Range("G66").Value = 0.08100650416169236
$ws.Range("H66").Value = 4.278046992874565
$ws.Range("G67").Value = 0.0810668810989057
$ws.Range("H67").Value = 2.784360022603085
$ws.Range("G68").Value = -0.01892728528162334
$ws.Range("H68").Value = 12.96076186840219
$ws.Range("G69").Value = -0.01879444970619375
$ws.Range("H69").Value = 1.812114481152396
$ws.Range("G70").Value = 0.08098689774141177
$ws.Range("H70").Value = 12.50389251182144
$ws.Range("G71").Value = 0.0623936121545007
$ws.Range("H71").Value = -21.43943991882328
$ws.Range("G72").Value = -0.143338017160853
$ws.Range("H72").Value = 6.706174032169021
$ws.Range("G73").Value = -0.1530682913110678
$ws.Range("H73").Value = 0.004745058374733753
$ws.Range("G74").Value = 0.1559921072651668
$ws.Range("H74").Value = 3.70492172689206
$ws.Range("G75").Value = 0.1547994515366237
$ws.Range("H75").Value = 2.892791221767383
$ws.Range("G76").Value = -0.01634907773615418
$ws.Range("H76").Value = -1477.393392605875
$ws.Range("G77").Value = -0.002590994360306139
$ws.Range("H77").Value = -17.35172312176479
$ws.Range("G78").Value = 0.09699144141080625
$ws.Range("H78").Value = 7.818283438902133
$ws.Range("G79").Value = 0.09491822157635314
$ws.Range("H79").Value = -2.046558533478845
$ws.Range("G80").Value = -0.222806072204085
$ws.Range("H80").Value = -2.944418599344977
$ws.Range("G81").Value = -0.2003133371677596
$ws.Range("H81").Value = 6.006435552716151
$ws.Range("G82").Value = 0.1710260554611701
$ws.Range("H82").Value = 2.033429585932346
$ws.Range("G83").Value = 0.1922974437183629
$ws.Range("H83").Value = 9.241953201532763
$ws.Range("G84").Value = 0.1151284590473831
$ws.Range("H84").Value = 8.479023577084217
$ws.Range("G85").Value = 0.1163039429050668
$ws.Range("H85").Value = 11.23304324899977

Write-Output "Applied all changes"